$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "Noun" will be created implicitly when used below (N2/N3).

# Column A: widen + apply a date/time number format to the whole column.
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(1).NumberFormat = "m/d/yy h:mm"

# Row 2
$ws.Range("A2").Value = 42609.642488425925
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 62
$ws.Range("D2").Value = 35
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 13710
$ws.Range("H2").Value = 10567
$ws.Range("I2").Value = 1697
$ws.Range("J2").Value = 201
$ws.Range("K2").Value = 116
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = "Noun"

# Row 3
$ws.Range("A3").Value = 42609.647094907406
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 62
$ws.Range("D3").Value = 36
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 13233
$ws.Range("H3").Value = 10546
$ws.Range("I3").Value = 1686
$ws.Range("J3").Value = 202
$ws.Range("K3").Value = 117
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = "Noun"
